$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("C2").Value = 2
$ws.Range("F2").Value = 55
$ws.Range("H2").Value = "living_rooms"
$ws.Range("L2").Value = "stimuli/img_lpr0l.png"
$ws.Range("M2").Value = 77.04651162790698
$ws.Range("N2").Value = 59.86046511627907
$ws.Range("O2").Value = 68.45348837209303
$ws.Range("P2").Value = 43
$ws.Range("Q2").Value = 7
$ws.Range("R2").Value = 7
$ws.Range("S2").Value = 7

# Row 3
$ws.Range("C3").Value = 2
$ws.Range("F3").Value = 56
$ws.Range("L3").Value = "stimuli/img_9jgbc.png"
$ws.Range("M3").Value = 40.30555555555556
$ws.Range("N3").Value = 24.08333333333333
$ws.Range("O3").Value = 32.19444444444444
$ws.Range("P3").Value = 36
$ws.Range("Q3").Value = 2
$ws.Range("R3").Value = 2
$ws.Range("S3").Value = 2

# Row 4
$ws.Range("C4").Value = 2
$ws.Range("F4").Value = 57
$ws.Range("H4").Value = "kitchens"
$ws.Range("L4").Value = "stimuli/img_6nbgt.png"
$ws.Range("M4").Value = 78.45161290322581
$ws.Range("N4").Value = 57.83870967741935
$ws.Range("O4").Value = 68.14516129032258
$ws.Range("P4").Value = 31

# Row 5
$ws.Range("C5").Value = 2
$ws.Range("F5").Value = 58
$ws.Range("H5").Value = "living_rooms"
$ws.Range("L5").Value = "stimuli/img_fmgjx.png"
$ws.Range("M5").Value = 79.90000000000001
$ws.Range("N5").Value = 56.975
$ws.Range("O5").Value = 68.4375
$ws.Range("P5").Value = 40
$ws.Range("Q5").Value = 7
$ws.Range("R5").Value = 7
$ws.Range("S5").Value = 7

# Row 6
$ws.Range("C6").Value = 2
$ws.Range("F6").Value = 59
$ws.Range("H6").Value = "kitchens"
$ws.Range("L6").Value = "stimuli/img_uy1n4.png"
$ws.Range("M6").Value = 76.30555555555556
$ws.Range("N6").Value = 55.33333333333334
$ws.Range("O6").Value = 65.81944444444444
$ws.Range("P6").Value = 36

# Row 7
$ws.Range("C7").Value = 2
$ws.Range("F7").Value = 60
$ws.Range("H7").Value = "kitchens"
$ws.Range("L7").Value = "stimuli/img_j5rpx.png"
$ws.Range("M7").Value = 72.24242424242425
$ws.Range("N7").Value = 50
$ws.Range("O7").Value = 61.12121212121212
$ws.Range("P7").Value = 33

# Row 8
$ws.Range("C8").Value = 2
$ws.Range("F8").Value = 61
$ws.Range("H8").Value = "kitchens"
$ws.Range("L8").Value = "stimuli/img_xti0z.png"
$ws.Range("M8").Value = 81.40625
$ws.Range("N8").Value = 61.4375
$ws.Range("O8").Value = 71.421875
$ws.Range("P8").Value = 32

# Row 9
$ws.Range("C9").Value = 2
$ws.Range("F9").Value = 62
$ws.Range("H9").Value = "kitchens"
$ws.Range("L9").Value = "stimuli/img_bwo9g.png"
$ws.Range("M9").Value = 64.81818181818181
$ws.Range("N9").Value = 42.36363636363637
$ws.Range("O9").Value = 53.59090909090909
$ws.Range("P9").Value = 33

# Row 10
$ws.Range("C10").Value = 2
$ws.Range("F10").Value = 63
$ws.Range("H10").Value = "kitchens"
$ws.Range("L10").Value = "stimuli/img_84s7n.png"
$ws.Range("M10").Value = 11.03125
$ws.Range("N10").Value = 2.90625
$ws.Range("O10").Value = 6.96875
$ws.Range("P10").Value = 32

# Row 11
$ws.Range("C11").Value = 2
$ws.Range("F11").Value = 64
$ws.Range("H11").Value = "kitchens"
$ws.Range("L11").Value = "stimuli/img_jz3kd.png"
$ws.Range("M11").Value = 72.79411764705883
$ws.Range("N11").Value = 51.64705882352941
$ws.Range("O11").Value = 62.22058823529412
$ws.Range("P11").Value = 34

# Row 12
$ws.Range("C12").Value = 2
$ws.Range("F12").Value = 65
$ws.Range("H12").Value = "living_rooms"
$ws.Range("L12").Value = "stimuli/img_89dvt.png"
$ws.Range("M12").Value = 81.09756097560975
$ws.Range("N12").Value = 64.6829268292683
$ws.Range("O12").Value = 72.89024390243902
$ws.Range("P12").Value = 41
$ws.Range("Q12").Value = 8
$ws.Range("R12").Value = 8
$ws.Range("S12").Value = 8

# Row 13
$ws.Range("C13").Value = 2
$ws.Range("F13").Value = 66
$ws.Range("H13").Value = "kitchens"
$ws.Range("L13").Value = "stimuli/img_p659z.png"
$ws.Range("M13").Value = 84.21621621621621
$ws.Range("N13").Value = 65.37837837837837
$ws.Range("O13").Value = 74.79729729729729
$ws.Range("P13").Value = 37
$ws.Range("R13").Value = 9
$ws.Range("S13").Value = 9

# Row 14
$ws.Range("C14").Value = 2
$ws.Range("F14").Value = 67
$ws.Range("H14").Value = "kitchens"
$ws.Range("L14").Value = "stimuli/img_05flq.png"
$ws.Range("M14").Value = 47.10344827586207
$ws.Range("N14").Value = 25.72413793103448
$ws.Range("O14").Value = 36.41379310344828
$ws.Range("P14").Value = 29
$ws.Range("R14").Value = 1
$ws.Range("S14").Value = 1

# Row 15
$ws.Range("C15").Value = 2
$ws.Range("F15").Value = 68
$ws.Range("H15").Value = "kitchens"
$ws.Range("L15").Value = "stimuli/img_i2k07.png"
$ws.Range("M15").Value = 64.25925925925925
$ws.Range("N15").Value = 40.92592592592592
$ws.Range("O15").Value = 52.59259259259259
$ws.Range("P15").Value = 27

# Row 16
$ws.Range("C16").Value = 2
$ws.Range("F16").Value = 69
$ws.Range("H16").Value = "kitchens"
$ws.Range("L16").Value = "stimuli/img_xesl0.png"
$ws.Range("M16").Value = 69.28571428571429
$ws.Range("N16").Value = 47.35714285714285
$ws.Range("O16").Value = 58.32142857142857
$ws.Range("P16").Value = 28

# Row 17
$ws.Range("C17").Value = 2
$ws.Range("F17").Value = 70
$ws.Range("H17").Value = "kitchens"
$ws.Range("L17").Value = "stimuli/img_c0me7.png"
$ws.Range("M17").Value = 68.40000000000001
$ws.Range("N17").Value = 45.62857142857143
$ws.Range("O17").Value = 57.01428571428572
$ws.Range("P17").Value = 35

# Row 18
$ws.Range("C18").Value = 2
$ws.Range("F18").Value = 71
$ws.Range("H18").Value = "bedrooms"
$ws.Range("L18").Value = "stimuli/img_ca8kd.png"
$ws.Range("M18").Value = 92.05405405405405
$ws.Range("N18").Value = 73.02702702702703
$ws.Range("O18").Value = 82.54054054054055
$ws.Range("P18").Value = 37
$ws.Range("Q18").Value = 10
$ws.Range("R18").Value = 10
$ws.Range("S18").Value = 10

# Row 19
$ws.Range("C19").Value = 2
$ws.Range("F19").Value = 72
$ws.Range("H19").Value = "kitchens"
$ws.Range("L19").Value = "stimuli/img_cv9qj.png"
$ws.Range("M19").Value = 60.34375
$ws.Range("N19").Value = 35.34375
$ws.Range("O19").Value = 47.84375
$ws.Range("P19").Value = 32

# Row 20
$ws.Range("C20").Value = 2
$ws.Range("F20").Value = 73
$ws.Range("H20").Value = "kitchens"
$ws.Range("L20").Value = "stimuli/img_ifebc.png"
$ws.Range("M20").Value = 84
$ws.Range("N20").Value = 65.88235294117646
$ws.Range("O20").Value = 74.94117647058823
$ws.Range("P20").Value = 34
$ws.Range("R20").Value = 9
$ws.Range("S20").Value = 9

# Row 21
$ws.Range("C21").Value = 2
$ws.Range("F21").Value = 74
$ws.Range("H21").Value = "kitchens"
$ws.Range("L21").Value = "stimuli/img_ncr40.png"
$ws.Range("M21").Value = 75.66666666666667
$ws.Range("N21").Value = 54.27272727272727
$ws.Range("O21").Value = 64.96969696969697
$ws.Range("P21").Value = 33

# Row 22
$ws.Range("C22").Value = 2
$ws.Range("F22").Value = 75
$ws.Range("H22").Value = "living_rooms"
$ws.Range("L22").Value = "stimuli/img_73pyk.png"
$ws.Range("M22").Value = 69.27659574468085
$ws.Range("N22").Value = 47.27659574468085
$ws.Range("O22").Value = 58.27659574468085
$ws.Range("Q22").Value = 5
$ws.Range("R22").Value = 5
$ws.Range("S22").Value = 5

# Row 23
$ws.Range("C23").Value = 2
$ws.Range("F23").Value = 76
$ws.Range("H23").Value = "kitchens"
$ws.Range("L23").Value = "stimuli/img_es7o2.png"
$ws.Range("M23").Value = 52.48571428571429
$ws.Range("N23").Value = 27.54285714285714
$ws.Range("O23").Value = 40.01428571428572
$ws.Range("P23").Value = 35

# Row 24
$ws.Range("C24").Value = 2
$ws.Range("F24").Value = 77
$ws.Range("H24").Value = "kitchens"
$ws.Range("L24").Value = "stimuli/img_ua9bs.png"
$ws.Range("M24").Value = 82
$ws.Range("N24").Value = 62.23333333333333
$ws.Range("O24").Value = 72.11666666666667
$ws.Range("P24").Value = 30

# Row 25
$ws.Range("C25").Value = 2
$ws.Range("F25").Value = 78
$ws.Range("H25").Value = "kitchens"
$ws.Range("L25").Value = "stimuli/img_c79r7.png"
$ws.Range("M25").Value = 56.26470588235294
$ws.Range("N25").Value = 34.26470588235294
$ws.Range("O25").Value = 45.26470588235294
$ws.Range("P25").Value = 34

# Row 26
$ws.Range("C26").Value = 2
$ws.Range("F26").Value = 79
$ws.Range("H26").Value = "kitchens"
$ws.Range("L26").Value = "stimuli/img_7pgd2.png"
$ws.Range("M26").Value = 78.59375
$ws.Range("N26").Value = 57.84375
$ws.Range("O26").Value = 68.21875
$ws.Range("P26").Value = 32
$ws.Range("R26").Value = 7
$ws.Range("S26").Value = 7

# Row 27
$ws.Range("C27").Value = 2
$ws.Range("F27").Value = 80
$ws.Range("H27").Value = "kitchens"
$ws.Range("L27").Value = "stimuli/img_s9are.png"
$ws.Range("M27").Value = 90.14285714285714
$ws.Range("N27").Value = 75.22857142857143
$ws.Range("O27").Value = 82.68571428571428
$ws.Range("P27").Value = 35
